$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 15 ("Willian Massami Watanabe"),
# pushing it down to row 17 and making room for the two new orientadores.
$ws.Rows.Item(15).Resize(2).Insert()

# Row 15: Reginaldo Fidelis
$ws.Range("A15").Value = "Reginaldo Fidelis"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 0

# Row 16: Rogerio Tondato
$ws.Range("A16").Value = "Rogerio Tondato"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 55
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 0
